$d = $word.ActiveDocument

$d.Content.Find.Execute("Seleksi Nasional Berbasis", $true, $false, $false, $false, $false, $true, 1, $false, "Seleksi Nasional Berdasarkan", 2)
